# Update "Date" (column B) values in row 2 of each regression-suite sheet to
# reflect the latest Multibill / CMCDeferred test execution timestamps, and
# flip the PayNow "No CF" (Corp-less card/PS) results from Fail to Pass.

$wb = $excel.ActiveWorkbook

$updates = @{
    "PayNowNoCFPC_27"   = @{ A2 = "Pass"; B2 = "Wed Mar 26 15:11:48 IST 2025" }
    "PayNowNoCFPS_27"   = @{ A2 = "Pass"; B2 = "Wed Mar 26 15:12:38 IST 2025" }
    "PayNowNoCFCorp_27" = @{ B2 = "Wed Mar 26 15:10:28 IST 2025" }
    "PayNowSCFPC_27"    = @{ B2 = "Wed Mar 26 15:15:52 IST 2025" }
    "PayNowSCFPS_27"    = @{ B2 = "Wed Mar 26 15:17:15 IST 2025" }
    "PayNowSCFCorp_27"  = @{ B2 = "Wed Mar 26 15:15:00 IST 2025" }
    "PayNowDCFPC_27"    = @{ B2 = "Wed Mar 26 15:06:09 IST 2025" }
    "PayNowDCFPS_27"    = @{ B2 = "Wed Mar 26 15:07:33 IST 2025" }
    "PayNowDCFCorp_27"  = @{ B2 = "Wed Mar 26 15:04:45 IST 2025" }
    "CCDeferredPS_27"   = @{ B2 = "Wed Mar 26 14:51:39 IST 2025" }
    "CCDeferredPC_27"   = @{ B2 = "Wed Mar 26 14:50:03 IST 2025" }
    "CCDeferredCorp_27" = @{ B2 = "Wed Mar 26 14:48:26 IST 2025" }
    "CMCAutopayPC_27"   = @{ B2 = "Wed Mar 26 14:58:29 IST 2025" }
    "CMCAutopayCorp_27" = @{ B2 = "Wed Mar 26 14:56:06 IST 2025" }
    "CMCAutopayPS_27"   = @{ B2 = "Wed Mar 26 15:00:51 IST 2025" }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellValues = $updates[$sheetName]
    foreach ($cellRef in $cellValues.Keys) {
        $ws.Range($cellRef).Value = $cellValues[$cellRef]
    }
}
